# Add a new "Swiss" worksheet (test data for the Switzerland market),
# cloned from the existing "Czech" sheet so it inherits the same layout,
# column widths, styles and merged cells, then update the market-specific
# cell values. Finally, clear the old "Germany" tab's selection/active
# state since the newly added sheet becomes the active tab.

$wb = $excel.ActiveWorkbook

# The previously active tab ("Germany") had a lingering cell selection
# and tabSelected flag; once the new sheet is activated, Excel drops
# tabSelected from Germany's view and resets its selection to the
# whole sheet (as if the user did "Select All" before switching away).
$germany = $wb.Worksheets.Item("Germany")
$germany.Activate()
$germany.Cells.Select()

# Clone the "Czech" sheet (same Wg/Miscellaneous/MPM800/PR1D2 template)
# to the end of the tab strip and rename it "Swiss".
$czech = $wb.Worksheets.Item("Czech")
$czech.Copy($null, $czech)
$swiss = $wb.Worksheets.Item("Czech (2)")
$swiss.Name = "Swiss"

# Fill in the Switzerland-specific market name and NGC/test-case ids.
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2653/T2656"

# Leave the cursor on B5 of the new sheet, matching the saved selection.
$swiss.Range("B5").Select()
